$d = $word.ActiveDocument

# Shared namespace / run-properties / paragraph-properties fragments used by
# both new "Paragraphe de liste" bullet items (same list + font formatting
# as the existing bullets in this document).
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$rPr = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>'
$pPr = '<w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' + $rPr + '</w:pPr>'

# New bullet 1: "Gérer l'absence d'un membre de l'équipe au sein du projet"
$para1 = '<w:p ' + $wNs + '>' + $pPr + `
    '<w:r>' + $rPr + '<w:t>Gérer l’absence d’un membre de l’équipe au sein du projet</w:t></w:r>' + `
    '</w:p>'

# New bullet 2: "Gérer le travail en distantiel au sein de l'équipe." — carries
# the relocated _GoBack bookmark (it was on the old last paragraph) and the
# spell-check proofErr markers around "distantiel" exactly as authored.
$para2 = '<w:p ' + $wNs + '>' + $pPr + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve">Gérer le travail en </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r>' + $rPr + '<w:t>distantiel</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> au sein de l’équipe.</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
    '</w:p>'

# Insert bullet 1 right before the final paragraph mark of the (current)
# last paragraph — this appends a new paragraph after it without disturbing
# that paragraph's own runs/attributes/bookmark.
$lastPara = $d.Paragraphs.Last
$ip = $d.Range($lastPara.Range.End - 1, $lastPara.Range.End - 1)
$ip.InsertXML($para1)

# Insert bullet 2 the same way, now after bullet 1.
$d2 = $word.ActiveDocument
$lastPara2 = $d2.Paragraphs.Last
$ip2 = $d2.Range($lastPara2.Range.End - 1, $lastPara2.Range.End - 1)
$ip2.InsertXML($para2)

# The _GoBack bookmark now exists twice (original spot + the one just added
# to bullet 2) — remove the original one so it only marks the new location.
$d3 = $word.ActiveDocument
$oldBookmark = $d3.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()
